$d = $word.ActiveDocument

$newText = "Die Schaubilder in diesem Dokument wurden von Jan Hollan, CzechGlobe, bereitgestellt. (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

# XML-escape the replacement text before embedding it in the InsertXML payload.
$newTextXml = $newText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Schaubilder*CzechGlobe*GaNight*") {

        # Range covering the run text only (exclude the trailing paragraph mark).
        $r = $p.Range
        $r.MoveEnd(1, -1)
        $oldLen = $r.End - $r.Start

        # Insert the replacement as raw WordprocessingML (an empty run followed
        # by a single unformatted run holding all the new text) right before the
        # existing runs - InsertXML merges the <w:r> children into the paragraph
        # that already exists at this position, leaving w:pPr untouched.
        $insertPoint = $d.Range($r.Start, $r.Start)
        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
               '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body><w:p><w:r/><w:r><w:t>' + $newTextXml + '</w:t></w:r></w:p></w:body>' +
               '</w:document></pkg:xmlData></pkg:part></pkg:package>'

        $insertPoint.InsertXML($xml)

        # The old runs got pushed right by the length of the newly inserted
        # text; delete them now that the replacement is in place.
        $newLen = $newText.Length
        $oldStart = $r.Start + $newLen
        $oldEnd = $oldStart + $oldLen
        $oldRange = $d.Range($oldStart, $oldEnd)
        $oldRange.Delete()

        break
    }
}
